$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the "Y" + "our " runs into a single "Your " run at the start of
#    the intro paragraph. MatchCase avoids touching the later, lower-case
#    "your" elsewhere in the same paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Your ", $true, $false, $false, $false, $false, $true, 1, $false, "Your ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Merge the tail runs ". " + "Below is ... your " + "estimated f" +
#    "ood " + "s" + "tamp amount." into a single run. The concatenation of
#    the existing runs is identical to the desired final text, so this is
#    purely a run-boundary collapse, not a content change.
#
#    Find's "Replace" argument runs Word's AutoCorrect/AutoFormat "smart
#    quotes" pass, which would turn the straight apostrophe in "doesn't"
#    into a curly one. To keep the original straight apostrophe, the
#    replacement is done in two steps: first collapse the runs using a
#    placeholder in place of the apostrophe, then patch the apostrophe back
#    in with a direct Range.Text assignment (which does not invoke
#    AutoCorrect).
# ---------------------------------------------------------------------------
$oldTail = ". Below is a detailed explanation of how this number was calculated. You can review it to make sure all the information is correct. If anything doesn't look right, you can return to the tool and correct or change an amount to see if it changes your estimated food stamp amount."
$placeholderTail = $oldTail.Replace("doesn't", "doesnZZZPLACEHOLDERZZZt")

$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $placeholderTail, 2) | Out-Null

$fix = $d.Content
$fix.Find.Execute("ZZZPLACEHOLDERZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fix.Text = "'"

# ---------------------------------------------------------------------------
# 3. Remove the old "_GoBack" bookmark that wrapped the "Calculation
#    Details" heading text.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 4. Add a new (collapsed) "_GoBack" bookmark immediately before the
#    "{{ paysAC " merge field, matching where Word last left the cursor.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("{{ paysAC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
